# Update conjoint positive figures after 85% data collected
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @{ B = 0.782258064516129;  C = 0.717305151915456;  D = 0.801104972375691;  E = 0.661490683229814;  F = 0.541808550889141 }
    3  = @{ B = 0.503597122302158;  C = 0.563106796116505;  D = 0.573643410852713;  E = 0.5;                 F = 0.491499227202473 }
    4  = @{ B = 0.855345911949686;  C = 0.735135135135135;  D = 0.824427480916031;  E = 0.734939759036145;  F = 0.537890044576523 }
    5  = @{ B = 0.771428571428571;  C = 0.752747252747253;  D = 0.783783783783784;  E = 0.738853503184713;  F = 0.550375939849624 }
    6  = @{ B = 0.884353741496599;  C = 0.809782608695652;  D = 0.859259259259259;  E = 0.865030674846626;  F = 0.609422492401216 }
    7  = @{ B = 0.44954128440367;   C = 0.599502487562189;  D = 0.707865168539326;  E = 0.480769230769231;  F = 0.53448275862069 }
    8  = @{ B = 0.57679180887372;   C = 0.552112676056338;  D = 0.710144927536232;  E = 0.5;                 F = 0.537051184110008 }
    9  = @{ B = 0.480645161290323;  C = 0.536327608982827;  D = 0.49171270718232;   E = 0.512422360248447;  F = 0.523148148148148 }
    10 = @{ B = 0.617741935483871;  C = 0.598414795244386;  D = 0.644567219152855;  E = 0.566770186335404;  F = 0.585905349794239 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
